$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "tata" column (C) for rows 13-43: TATA 2 -> TATA 1
$ws.Range("C13:C43").Value = "TATA 1"

# Prenom_Nom_Promoteur (G) for rows 19-23: "Autre" -> "Mariama SADIO"
$ws.Range("G19:G23").Value = "Mariama SADIO"

# Precisez (H) for rows 19-23: clear the free-text detail now that the name is selected directly
$ws.Range("H19:H23").ClearContents()

# Reflect the cell the user ended up on when saving
$ws.Range("G15").Select()
